# ComplianceFormTemplate.docx edit:
#  1) Remove the "_GoBack" bookmark pair from the paragraph right after the
#     "DESCRIPTION OF FINDINGS" table (it keeps its pPr, but becomes empty).
#  2) Right before the final <w:sectPr>, insert a brand-new empty paragraph
#     (no pPr) ahead of the existing trailing paragraph, and turn that
#     trailing paragraph (which used to just hold a tab-stop pPr) into a
#     bare paragraph that now carries the "_GoBack" bookmark pair instead.

$d = $word.ActiveDocument

# --- Step 1: drop the bookmark from its original location -------------
# Deleting it first frees up bookmark id "0" so the re-inserted bookmark
# below naturally gets the same id back.
$bm = $d.Bookmarks.Item("_GoBack")
[void]$bm.Delete()

# --- Step 2: insert a brand-new bare paragraph before the last paragraph
$lastCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastCount)
[void]$lastPara.Range.InsertParagraphBefore()

# The freshly-inserted paragraph inherited the old tab-stop formatting;
# overwrite it in place with a genuinely bare <w:p/>.
$newEmptyIndex = $d.Paragraphs.Count - 1
$newEmpty = $d.Paragraphs.Item($newEmptyIndex)
$bareParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
[void]$newEmpty.Range.InsertXML($bareParagraphXml)

# --- Step 3: replace the old trailing (tab-stop) paragraph with one that
# only contains the "_GoBack" bookmark, and no pPr at all.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$trailing.Range.InsertXML($bookmarkParagraphXml)
